$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.790.90"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "1.784.38"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.18"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5121"
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3853"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07818"
$ws.Range("E9").Value = "  -8.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.089"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.58"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.202"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.14"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "1.775.20"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.213"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.52"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("E18").Value = "  -5.85%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").Value = "27.835.81"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.93"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "1.987.42"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.366"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.54"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1074"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.036"
$ws.Range("E32").Value = "  -5.65%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07087"
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02307"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.784"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2122"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.47"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6089"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.154"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("E44").Value = "  -5.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.09"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5891"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.704"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.99"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.201"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.896"
$ws.Range("E50").Value = "  -4.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06848"
$ws.Range("E51").Value = "  -1.91%  "
